# Journal de travail - modification
$wb = $excel.ActiveWorkbook

$wsJournal = $wb.Worksheets.Item("Journal")
$wsTotaux  = $wb.Worksheets.Item("Totaux")

# Correct the time logged for the Documentation entry on week 2 (Journal!C6): 3h instead of 2h15
$wsJournal.Range("C6").Value = 0.125

# Add the missing weekly total for week 2 in Totaux!B3
$wsTotaux.Range("B3").Formula = "=SUM(Journal!C6:C9)"

[void]$wsTotaux.Range("C8").Select()
[void]$wsJournal.Select()
[void]$wsJournal.Range("C18").Select()

$wb.Application.Calculate() | Out-Null
